$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = -7.824999999999998
$ws.Range("D8").Value = -8.048999999999999
$ws.Range("C12").Value = -12.977
$ws.Range("D12").Value = -7.913999999999999
$ws.Range("D14").Value = -8.263
$ws.Range("D22").Value = -8.191999999999998
